# Generate Report for Handback
# Adds a new handback record (0498cf13-025e-456d-aaa9-9f7d56f20919.md) to the
# "Overview", "zh-cn" and "de-de" worksheets/tables.

$wb = $excel.ActiveWorkbook

$fileGuid   = "0498cf13-025e-456d-aaa9-9f7d56f20919"
$fileName   = "$fileGuid.md"
$pathName   = "e2e\$fileGuid.md"
$statusSync = "Handed back: in sync with en-US"

$srcRepoBase   = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7aa8312b3bd6b607465aabbf87977f3f3c189c92/e2e/$fileName"
$zhcnRepoBase  = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/c9b6beee346c85c44b3b3e38d53827ace38d1daf/e2e/$fileName"
$dedeRepoBase  = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/c9b6beee346c85c44b3b3e38d53827ace38d1daf/e2e/$fileName"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$rowOverview = $loOverview.ListRows.Add()
$rngOverview = $rowOverview.Range

$rngOverview.Cells.Item(1, 1).Value = $fileName
$rngOverview.Cells.Item(1, 2).Value = $pathName
$rngOverview.Cells.Item(1, 3).Value = ".md"
$rngOverview.Cells.Item(1, 5).Value = $statusSync
$rngOverview.Cells.Item(1, 6).Value = $statusSync
$rngOverview.Cells.Item(1, 7).Value = "2016-09-07 09:56:31"

$wsOverview.Hyperlinks.Add($rngOverview.Cells.Item(1, 2), $srcRepoBase, "", "", $pathName) | Out-Null

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$rowZhCn = $loZhCn.ListRows.Add()
$rngZhCn = $rowZhCn.Range

$zhCnXlf = "$fileGuid.c9b6beee346c85c44b3b3e38d53827ace38d1daf.zh-cn.xlf"

$rngZhCn.Cells.Item(1, 1).Value = $fileName
$rngZhCn.Cells.Item(1, 2).Value = ".md"
$rngZhCn.Cells.Item(1, 3).Value = $statusSync
$rngZhCn.Cells.Item(1, 4).Value = "e2e"
$rngZhCn.Cells.Item(1, 5).Value = "ht"
$rngZhCn.Cells.Item(1, 6).Value = "True"
$rngZhCn.Cells.Item(1, 7).Value = $zhCnXlf
$rngZhCn.Cells.Item(1, 8).Value = "2016-09-07 09:56:21"
$rngZhCn.Cells.Item(1, 9).Value = $fileName
$rngZhCn.Cells.Item(1, 10).Value = $zhCnXlf
$rngZhCn.Cells.Item(1, 11).Value = "2016-09-07 09:57:08"
$rngZhCn.Cells.Item(1, 12).Value = ""
$rngZhCn.Cells.Item(1, 13).Value = "True"
$rngZhCn.Cells.Item(1, 14).Value = ""
$rngZhCn.Cells.Item(1, 15).Value = "False"
$rngZhCn.Cells.Item(1, 16).Value = ""

$wsZhCn.Hyperlinks.Add($rngZhCn.Cells.Item(1, 1), $srcRepoBase, "", "", $fileName) | Out-Null
$wsZhCn.Hyperlinks.Add($rngZhCn.Cells.Item(1, 9), $zhcnRepoBase, "", "", $fileName) | Out-Null

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$rowDeDe = $loDeDe.ListRows.Add()
$rngDeDe = $rowDeDe.Range

$deDeXlf = "$fileGuid.c9b6beee346c85c44b3b3e38d53827ace38d1daf.de-de.xlf"

$rngDeDe.Cells.Item(1, 1).Value = $fileName
$rngDeDe.Cells.Item(1, 2).Value = ".md"
$rngDeDe.Cells.Item(1, 3).Value = $statusSync
$rngDeDe.Cells.Item(1, 4).Value = "e2e"
$rngDeDe.Cells.Item(1, 5).Value = "ht"
$rngDeDe.Cells.Item(1, 6).Value = "True"
$rngDeDe.Cells.Item(1, 7).Value = $deDeXlf
$rngDeDe.Cells.Item(1, 8).Value = "2016-09-07 09:56:31"
$rngDeDe.Cells.Item(1, 9).Value = $fileName
$rngDeDe.Cells.Item(1, 10).Value = $deDeXlf
$rngDeDe.Cells.Item(1, 11).Value = "2016-09-07 09:57:30"
$rngDeDe.Cells.Item(1, 12).Value = ""
$rngDeDe.Cells.Item(1, 13).Value = "True"
$rngDeDe.Cells.Item(1, 14).Value = ""
$rngDeDe.Cells.Item(1, 15).Value = "False"
$rngDeDe.Cells.Item(1, 16).Value = ""

$wsDeDe.Hyperlinks.Add($rngDeDe.Cells.Item(1, 1), $srcRepoBase, "", "", $fileName) | Out-Null
$wsDeDe.Hyperlinks.Add($rngDeDe.Cells.Item(1, 9), $dedeRepoBase, "", "", $fileName) | Out-Null
